$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.346.84"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "3.101.51"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.46"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.64"
$ws.Range("E6").Value = "  -1.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.098.91"
$ws.Range("E8").Value = "  -1.16%  "
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("E11").Value = "  -2.45%  "
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("E13").Value = "  -1.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.63"
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("D16").Value = "3.616.87"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").Value = "67.245.43"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.11"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.68"
$ws.Range("E19").Value = "  +3.22%  "
$ws.Range("D20").Value = "3.097.69"
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "489.50"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.81"
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("E23").Value = "  -2.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.64"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("E25").Value = "  -1.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").Value = "  -2.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.48"
$ws.Range("E27").Value = "  +4.55%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  -2.51%  "
$ws.Range("E30").Value = "  -3.17%  "
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.27"
$ws.Range("E32").Value = "  -2.83%  "
$ws.Range("E33").Value = "  -1.54%  "
$ws.Range("D34").Value = "0.0₃0942"
$ws.Range("E34").Value = "  -5.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  -2.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.965"
$ws.Range("E37").Value = "  -2.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "46.29"
$ws.Range("E38").Value = "  -2.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.02"
$ws.Range("E39").Value = "  -4.34%  "
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("E41").Value = "  -2.20%  "
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("D43").Value = "2.795.02"
$ws.Range("E43").Value = "  -1.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "380.21"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.58"
$ws.Range("E45").Value = "  -7.07%  "
$ws.Range("E46").Value = "  -2.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "135.69"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.80"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("E51").Value = "  -1.58%  "
